$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.183.53'
$ws.Range("E2").Value = '  -0.97%  '

$ws.Range("D3").Value = '1.869.17'
$ws.Range("E3").Value = '  -2.24%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''233.77'
$ws.Range("E5").Value = '  -2.37%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = '''0.4683'
$ws.Range("E7").Value = '  -2.23%  '

$ws.Range("D8").Value = '''0.2830'
$ws.Range("E8").Value = '  -0.40%  '

$ws.Range("D9").Value = '''0.06586'
$ws.Range("E9").Value = '  -1.64%  '

$ws.Range("D10").Value = '''20.25'
$ws.Range("E10").Value = '  +7.35%  '

$ws.Range("D11").Value = '''0.07780'
$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("D12").Value = '''96.90'
$ws.Range("E12").Value = '  -5.18%  '

$ws.Range("D13").Value = '1.880.28'
$ws.Range("E13").Value = '  -1.70%  '

$ws.Range("D14").Value = '''5.071'
$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").Value = '''0.6693'

$ws.Range("D16").Value = '''286.16'
$ws.Range("E16").Value = '  +5.37%  '

$ws.Range("D17").Value = '30.220.05'
$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("D18").Value = '''1.001'
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").Value = '2.133.53'
$ws.Range("E19").Value = '  -1.17%  '

$ws.Range("D20").Value = '''12.56'

$ws.Range("D21").Value = '''5.371'
$ws.Range("E21").Value = '  -0.81%  '

$ws.Range("D22").Value = '''0.000007245'
$ws.Range("E22").Value = '  -3.03%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").Value = '''6.155'
$ws.Range("E24").Value = '  -2.31%  '

$ws.Range("D25").Value = '''9.335'
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("D26").Value = '''167.31'
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").Value = '''19.09'
$ws.Range("E27").Value = '  -0.80%  '

$ws.Range("D28").Value = '''1.959'
$ws.Range("E28").Value = '  -4.81%  '

$ws.Range("E29").Value = '  -1.41%  '

$ws.Range("D30").Value = '''0.09652'
$ws.Range("E30").Value = '  -3.61%  '

$ws.Range("D31").Value = '''4.382'
$ws.Range("E31").Value = '  -5.02%  '

$ws.Range("D32").Value = '''1.465'
$ws.Range("E32").Value = '  -3.29%  '

$ws.Range("D33").Value = '''4.107'
$ws.Range("E33").Value = '  -2.81%  '

$ws.Range("D34").Value = '''0.04668'
$ws.Range("E34").Value = '  -1.24%  '

$ws.Range("D35").Value = '''0.7025'
$ws.Range("E35").Value = '  -3.24%  '

$ws.Range("D36").Value = '''1.088'
$ws.Range("E36").Value = '  -1.69%  '

$ws.Range("D37").Value = '''1.001'
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").Value = '''2.720'
$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("D39").Value = '''0.01864'
$ws.Range("E39").Value = '  -2.57%  '

$ws.Range("D40").Value = '''6.429'
$ws.Range("E40").Value = '  +2.05%  '

$ws.Range("D41").Value = '''2.521'
$ws.Range("E41").Value = '  -3.48%  '

$ws.Range("D42").Value = '''71.57'
$ws.Range("E42").Value = '  -4.31%  '

$ws.Range("D43").Value = '''0.8586'
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").Value = '''1.943'
$ws.Range("E44").Value = '  -1.24%  '

$ws.Range("D45").Value = '''1.001'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").Value = '''102.99'
$ws.Range("E46").Value = '  -2.06%  '

$ws.Range("D47").Value = '''0.4179'
$ws.Range("E47").Value = '  -1.99%  '

$ws.Range("D48").Value = '''987.11'
$ws.Range("E48").Value = '  +7.71%  '

$ws.Range("D49").Value = '''7.178'
$ws.Range("E49").Value = '  -3.02%  '

$ws.Range("D50").Value = '''9.170'
$ws.Range("E50").Value = '  +4.78%  '

$ws.Range("D51").Value = '''33.88'
$ws.Range("E51").Value = '  -2.57%  '
